# Crabs for RNA analysis - normalize "Exposure duration" labels to lowercase
# (column C, "Long" -> "long", "Short" -> "short") across all data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C holds the "Exposure duration" values (rows 9-71). Replace the
# capitalized "Long"/"Short" labels with their lowercase equivalents.
$rng = $ws.Range("C9:C71")
$rng.Replace("Long", "long")
$rng.Replace("Short", "short")

# Reflect the author's final selection/scroll position from the edit session.
$ws.Range("H31").Select()
